# Update forecast error table values for rows 7-11 (B:G) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    7  = @(0.188856333554674, 0.6210211877426164, 0.8227230598494827, 0.9070408259000708, 0.8997464472128907, 36)
    8  = @(0.2018765014077513, 0.5836206403444438, 0.7462698502533404, 0.8638691163905214, 0.8522125269166538, 35)
    9  = @(-0.09022918269035383, 0.3486329789925496, 0.2017908439045041, 0.4492113577198423, 0.451488227386921, 20)
    10 = @(-0.009659961070461246, 0.414243314843904, 0.2785379243470545, 0.5277669223691975, 0.5492252060470605, 13)
    11 = @(0.02760926664935082, 0.3666891604559107, 0.2015955073104189, 0.4489938833775121, 0.5010404620705597, 5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    # Columns B(2) through G(7)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

$wb.Save()
